$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$detailQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
 WHERE ss.disease_subtype IN ["Medullary Carcinoma"] and sf.grouped_recurrence_score IN ["16-20"] and d.tumor_size_group In ["(2,3]"] and tp.chemotherapy_regimen In ["IV CMF (3 week cycles)"]
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`
'@

$webCountQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Medullary Carcinoma"] and sf.grouped_recurrence_score IN ["16-20"] and d.tumor_size_group In ["(2,3]"] and tp.chemotherapy_regimen In ["IV CMF (3 week cycles)"]
WITH ss
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
'@

$newCountQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Medullary Carcinoma"] and sf.grouped_recurrence_score IN ["16-20"]and d.tumor_size_group In ["(2,3]"] and tp.chemotherapy_regimen In ["IV CMF (3 week cycles)"]
WITH ss
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
'@

$ws.Range("B2").Value2 = $detailQuery
$ws.Range("C2").Value2 = $webCountQuery
$ws.Range("B3").Value2 = $detailQuery
$ws.Range("C3").Value2 = $newCountQuery
$ws.Range("B4").Value2 = $detailQuery
$ws.Range("C4").Value2 = $webCountQuery

$ws.Range("C4").Select() | Out-Null
